$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Cells.Item(2, 4) "30.352.96"
Set-TextCell $ws.Cells.Item(2, 5) "  +0.77%  "

# Row 3
Set-TextCell $ws.Cells.Item(3, 4) "1.934.88"
Set-TextCell $ws.Cells.Item(3, 5) "  +1.09%  "

# Row 4
Set-TextCell $ws.Cells.Item(4, 5) "  -0.08%  "

# Row 5
Set-TextCell $ws.Cells.Item(5, 4) "251.79"
Set-TextCell $ws.Cells.Item(5, 5) "  +2.64%  "

# Row 6
Set-TextCell $ws.Cells.Item(6, 4) "0.7226"
Set-TextCell $ws.Cells.Item(6, 5) "  +3.10%  "

# Row 8
Set-TextCell $ws.Cells.Item(8, 4) "0.3296"
Set-TextCell $ws.Cells.Item(8, 5) "  +2.44%  "

# Row 9
Set-TextCell $ws.Cells.Item(9, 4) "27.78"
Set-TextCell $ws.Cells.Item(9, 5) "  +7.51%  "

# Row 10
Set-TextCell $ws.Cells.Item(10, 4) "0.07231"
Set-TextCell $ws.Cells.Item(10, 5) "  +5.66%  "

# Row 11
Set-TextCell $ws.Cells.Item(11, 4) "0.8056"
Set-TextCell $ws.Cells.Item(11, 5) "  +2.34%  "

# Row 12
Set-TextCell $ws.Cells.Item(12, 4) "0.08092"
Set-TextCell $ws.Cells.Item(12, 5) "  +2.04%  "

# Row 13
Set-TextCell $ws.Cells.Item(13, 4) "1.934.95"
Set-TextCell $ws.Cells.Item(13, 5) "  +1.07%  "

# Row 14
Set-TextCell $ws.Cells.Item(14, 4) "5.456"
Set-TextCell $ws.Cells.Item(14, 5) "  +1.87%  "

# Row 15
Set-TextCell $ws.Cells.Item(15, 4) "94.71"
Set-TextCell $ws.Cells.Item(15, 5) "  +1.35%  "

# Row 16
Set-TextCell $ws.Cells.Item(16, 4) "15.06"
Set-TextCell $ws.Cells.Item(16, 5) "  +5.22%  "

# Row 17
Set-TextCell $ws.Cells.Item(17, 4) "30.350.37"
Set-TextCell $ws.Cells.Item(17, 5) "  +0.69%  "

# Row 18
Set-TextCell $ws.Cells.Item(18, 2) "BitcoinCash"
Set-TextCell $ws.Cells.Item(18, 3) "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell $ws.Cells.Item(18, 4) "253.71"
Set-TextCell $ws.Cells.Item(18, 5) "  -2.25%  "

# Row 19
Set-TextCell $ws.Cells.Item(19, 2) "ShibaInu"
Set-TextCell $ws.Cells.Item(19, 3) "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell $ws.Cells.Item(19, 4) "0.000008225"
Set-TextCell $ws.Cells.Item(19, 5) "  +5.03%  "

# Row 20
Set-TextCell $ws.Cells.Item(20, 4) "5.824"
Set-TextCell $ws.Cells.Item(20, 5) "  +0.72%  "

# Row 21
Set-TextCell $ws.Cells.Item(21, 4) "2.189.91"
Set-TextCell $ws.Cells.Item(21, 5) "  +0.89%  "

# Row 22
Set-TextCell $ws.Cells.Item(22, 5) "  -0.01%  "

# Row 23
Set-TextCell $ws.Cells.Item(23, 5) "  -0.04%  "

# Row 24
Set-TextCell $ws.Cells.Item(24, 4) "6.947"
Set-TextCell $ws.Cells.Item(24, 5) "  +1.99%  "

# Row 25
Set-TextCell $ws.Cells.Item(25, 4) "9.727"

# Row 26
Set-TextCell $ws.Cells.Item(26, 4) "166.09"
Set-TextCell $ws.Cells.Item(26, 5) "  +3.66%  "

# Row 27
Set-TextCell $ws.Cells.Item(27, 4) "2.346"
Set-TextCell $ws.Cells.Item(27, 5) "  +5.65%  "

# Row 28
Set-TextCell $ws.Cells.Item(28, 4) "19.30"
Set-TextCell $ws.Cells.Item(28, 5) "  +3.21%  "

# Row 29
Set-TextCell $ws.Cells.Item(29, 4) "0.1295"
Set-TextCell $ws.Cells.Item(29, 5) "  -0.87%  "

# Row 30
Set-TextCell $ws.Cells.Item(30, 4) "1.356"
Set-TextCell $ws.Cells.Item(30, 5) "  -0.16%  "

# Row 31
Set-TextCell $ws.Cells.Item(31, 4) "1.545"
Set-TextCell $ws.Cells.Item(31, 5) "  -0.23%  "

# Row 32
Set-TextCell $ws.Cells.Item(32, 4) "4.444"
Set-TextCell $ws.Cells.Item(32, 5) "  +1.26%  "

# Row 33
Set-TextCell $ws.Cells.Item(33, 5) "  +0.82%  "

# Row 34
Set-TextCell $ws.Cells.Item(34, 4) "0.05242"
Set-TextCell $ws.Cells.Item(34, 5) "  +4.36%  "

# Row 35
Set-TextCell $ws.Cells.Item(35, 5) "  +6.53%  "

# Row 36
Set-TextCell $ws.Cells.Item(36, 4) "0.7512"
Set-TextCell $ws.Cells.Item(36, 5) "  +1.58%  "

# Row 37
Set-TextCell $ws.Cells.Item(37, 4) "2.760"
Set-TextCell $ws.Cells.Item(37, 5) "  +1.31%  "

# Row 38
Set-TextCell $ws.Cells.Item(38, 4) "0.01969"
Set-TextCell $ws.Cells.Item(38, 5) "  +3.11%  "

# Row 39
Set-TextCell $ws.Cells.Item(39, 4) "2.800"
Set-TextCell $ws.Cells.Item(39, 5) "  +0.54%  "

# Row 40
Set-TextCell $ws.Cells.Item(40, 4) "79.22"
Set-TextCell $ws.Cells.Item(40, 5) "  -0.04%  "

# Row 41
Set-TextCell $ws.Cells.Item(41, 4) "6.444"
Set-TextCell $ws.Cells.Item(41, 5) "  -0.78%  "

# Row 42
Set-TextCell $ws.Cells.Item(42, 4) "0.4541"
Set-TextCell $ws.Cells.Item(42, 5) "  +3.19%  "

# Row 43
Set-TextCell $ws.Cells.Item(43, 4) "2.031"
Set-TextCell $ws.Cells.Item(43, 5) "  +1.42%  "

# Row 44
Set-TextCell $ws.Cells.Item(44, 4) "0.8426"
Set-TextCell $ws.Cells.Item(44, 5) "  +1.40%  "

# Row 45
Set-TextCell $ws.Cells.Item(45, 5) "  -0.02%  "

# Row 46
Set-TextCell $ws.Cells.Item(46, 4) "102.05"
Set-TextCell $ws.Cells.Item(46, 5) "  +0.47%  "

# Row 47
Set-TextCell $ws.Cells.Item(47, 4) "9.850"
Set-TextCell $ws.Cells.Item(47, 5) "  +1.68%  "

# Row 48
Set-TextCell $ws.Cells.Item(48, 4) "7.463"
Set-TextCell $ws.Cells.Item(48, 5) "  +3.92%  "

# Row 49
Set-TextCell $ws.Cells.Item(49, 4) "36.82"
Set-TextCell $ws.Cells.Item(49, 5) "  +2.84%  "

# Row 50
Set-TextCell $ws.Cells.Item(50, 4) "0.4192"
Set-TextCell $ws.Cells.Item(50, 5) "  +3.66%  "

# Row 51
Set-TextCell $ws.Cells.Item(51, 4) "0.06053"
Set-TextCell $ws.Cells.Item(51, 5) "  +2.43%  "
